# Add a "contratante" column to the database sheet and tag every
# existing record with the contracting company name ("Puma").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last populated row in column A (16 rows of data in this file).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Header for the new column.
$ws.Cells.Item(1, 3).Value = "contratante"

# Fill every data row with the contratante value.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = "Puma"
}

$ws.Range("C2:C16").Select()
